$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark additional rows as "Completed" in the Status column (F) ---
$ws.Range("F7").Value = "Completed"
$ws.Range("F21").Value = "Completed"
$ws.Range("F23").Value = "Completed"
$ws.Range("F26").Value = "Completed"
$ws.Range("F27").Value = "Completed"

# --- Rework the "powerups" requirement row (row 28) ---
$ws.Range("B28").Value = "Super shot power"
$ws.Range("C28").Value = 7

# --- Preserve the old row 29 ("Decide on theme") by copying it down to row 32,
#     which is where it lands in the restructured sheet. Do this BEFORE we
#     overwrite row 29 with new content below. ---
$ws.Range("A29:F29").Copy($ws.Range("A32:F32"))

# --- Replace row 29 with the new "Triple shot power" requirement.
#     Clear the old formatting/formula first so the row matches a freshly
#     typed (unstyled) row. ---
$ws.Range("A29:F29").Clear()
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "Triple shot power"
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 20
$ws.Range("E29").Value = 4
$ws.Rows("29").RowHeight = 15.75

# --- New row 30: "Speed up power" ---
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "Speed up power"
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 20
$ws.Range("E30").Value = 4
$ws.Rows("30").RowHeight = 15.75

# --- New row 31: "Overall Polishing" ---
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Overall Polishing"
$ws.Range("C31").Value = 7
$ws.Range("D31").Value = 10
$ws.Range("E31").Value = 4
$ws.Rows("31").RowHeight = 15.75

# --- Row 32 ("Decide on theme") now holds a literal id (no longer part of
#     the shared A-column formula series), since it moved out of the
#     contiguous numbered block. ---
$ws.Range("A32").Value = 29

# --- Move the "Total story points estimated:" summary row from 32 down to
#     34, leaving row 33 blank, and recompute the sum range. ---
$ws.Range("B32").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("B34").Value = "Total story points estimated:"
$ws.Range("C34").Formula = "=SUM(C3:C32)"

# --- Update the active selection/view to match the final state. ---
$ws.Range("F30").Select()
